$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.15874635842922658
$ws.Range("A2").Value = -0.047301588604405254
$ws.Range("A3").Value = -0.049544986139475
$ws.Range("A4").Value = 0.28398814138762063
$ws.Range("A5").Value = -0.0059999991275603293
$ws.Range("A6").Value = -0.0059999990921539847
$ws.Range("A7").Value = -0.019999998943845299
$ws.Range("A8").Value = -0.019999998935734453
$ws.Range("A9").Value = -0.0059999990726691266
$ws.Range("A10").Value = -0.0059999990657573221
$ws.Range("A11").Value = -0.0044999990804477363
$ws.Range("A12").Value = -0.005999999064228323
$ws.Range("A13").Value = -0.0059999990595578367
$ws.Range("A14").Value = -0.011999998996443217
$ws.Range("A15").Value = -0.0059999990587868979
$ws.Range("A16").Value = 0.023367834197445703
$ws.Range("A17").Value = -0.0059999990558816663
$ws.Range("A18").Value = -0.008999999023942884
$ws.Range("A19").Value = -0.032484823775399629
$ws.Range("A20").Value = -0.0089999991013858249
$ws.Range("A21").Value = -0.0089999991003510971
$ws.Range("A22").Value = -0.008999999099604139
$ws.Range("A23").Value = -0.0089999990966127541
$ws.Range("A24").Value = -0.04199999873939575
$ws.Range("A25").Value = -0.041999998732640265
$ws.Range("A26").Value = -0.005999999089755903
$ws.Range("A27").Value = -0.0059999990849641804
$ws.Range("A28").Value = -0.0059999990637669143
$ws.Range("A29").Value = -0.011999998987850091
$ws.Range("A30").Value = -0.019999998898175164
$ws.Range("A31").Value = -0.014999998940558257
$ws.Range("A32").Value = -0.020999998877056392
$ws.Range("A33").Value = -0.0059999990322632257
